$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New USCDI requirement rows appended to the table (rows 49-56).
$newRows = @(
    @("True", "Care Plan", "US Core CarePlan Profile", "CarePlan.addresses", "US Core CarePlan Profile-CarePlan.addresses"),
    @("True", "Family Health History", "US Core FamilyMemberHistory Profile", "FamilyMemberHistory.extension:recorder", "US Core FamilyMemberHistory Profile-FamilyMemberHistory.extension:recorder"),
    @("True", "Unique Device Identifier", "US Core Device Profile", "Device.udiCarrier", "US Core Device Profile-Device.udiCarrier"),
    @("True", "Unique Device Identifier", "US Core Device Profile", "Device.distinctIdentifier", "US Core Device Profile-Device.distinctIdentifier"),
    @("True", "Unique Device Identifier", "US Core Device Profile", "Device.manufactureDate", "US Core Device Profile-Device.manufactureDate"),
    @("True", "Unique Device Identifier", "US Core Device Profile", "Device.expirationDate", "US Core Device Profile-Device.expirationDate"),
    @("True", "Unique Device Identifier", "US Core Device Profile", "Device.lotNumber", "US Core Device Profile-Device.lotNumber"),
    @("True", "Unique Device Identifier", "US Core Device Profile", "Device.serialNumber", "US Core Device Profile-Device.serialNumber")
)

$startRow = 49
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    # Column A holds the literal text "True" (matches the rest of the sheet's
    # Is_New column, which stores plain strings rather than booleans). A
    # leading apostrophe forces Excel to store it as text instead of coercing
    # it into a native Boolean; resetting the style afterwards keeps the
    # quote-prefix flag from leaking into the cell's formatting.
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value = "'" + $row[0]
    $aCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
